$d = $word.ActiveDocument

# --- 1) "Wiplach" paragraph: drop the spell-check proofErr wrapper around the run ---
# A plain in-place text replace on a single fully-wrapped run leaves the (now orphaned)
# w:proofErr start/end markers behind, so instead we splice in a brand new paragraph
# holding the same text and then remove the old (proofErr-wrapped) paragraph outright.
$startRng = $d.Range(0, 0)
$startRng.InsertParagraphBefore()
$d.Paragraphs(1).Range.InsertBefore("Wiplach")
$d.Paragraphs(2).Range.Delete()

# --- 2) "Toy story 1" paragraph: merge the 3 runs ("Toy " / "story" / " 1") back into a
#        single run and drop the w:proofErr wrap around "story" ---
$found = $d.Content.Find.Execute("Toy story 1", $false, $false, $false, $false, $false, `
                                  $true, 1, $false, "Toy story 1", 2)

# --- 3) Insert a brand new paragraph "Toy story 2" right after "Toy story 1" and before
#        the trailing blank paragraph ---
$p2 = $d.Paragraphs(2)
$endOfP2 = $p2.Range
$endOfP2.Collapse(0)
$endOfP2.InsertParagraphAfter()
$d.Paragraphs(3).Range.InsertBefore("Toy story 2")

Write-Output "done"
